$d = $word.ActiveDocument

# --- Change 1: "...important task to that employee..." -> "...important tasks to that employee..."
# Insert the letter "s" right after "task" (before " to that employee") as a tracked
# insertion so Word naturally splits the surrounding run the way a genuine edit would,
# then accept just that one revision so the final OOXML ends up with three plain runs:
# "...important task", "s", " to that employee...that way." (the last one carrying
# xml:space="preserve" because it begins with a space).
$d.TrackRevisions = $true

$r1 = $d.Content
$r1.Find.Execute("important task to that employee", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$splitPos = $r1.Start + ("important task").Length
$insPoint = $d.Range($splitPos, $splitPos)
$insPoint.InsertAfter("s")

$d.TrackRevisions = $false
$d.Revisions(1).Accept()

# --- Change 2: fix typo "previously employer" -> "previous employer"
$d.Content.Find.Execute("previously employer", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "previous employer", 2)
